$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, shifting the existing rows 152-249 down to 153-250
$ws.Rows("152:152").Insert()

# Populate the newly inserted row 152 with the new price-report record
$ws.Range("A152").Value2 = 4
$ws.Range("B152").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value2 = "Los Lagos"
$ws.Range("D152").Value2 = 44762
$ws.Range("E152").Value2 = 10
$ws.Range("F152").Value2 = "Fruta"
$ws.Range("G152").Value2 = 100102
$ws.Range("H152").Value2 = "Cítricos"
$ws.Range("I152").Value2 = 100102004
$ws.Range("J152").Value2 = "Mandarina"
$ws.Range("K152").Value2 = "Clemenuless"
$ws.Range("L152").Value2 = "Primera"
$ws.Range("M152").Value2 = 300
$ws.Range("N152").Value2 = 8000
$ws.Range("O152").Value2 = 9000
$ws.Range("P152").Value2 = 8500
$ws.Range("Q152").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R152").Value2 = "Provincia de Limarí"
$ws.Range("S152").Value2 = 850
$ws.Range("T152").Value2 = 10
